$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated SOAE N_xi fitted parameter results after fixing Tokay fs error.
# Data for rows 2-49 (Anole/Owl/Tokay blocks) is replaced in full:
#  - Tokay block (regenerated values) now occupies rows 2-17
#  - Anole block (unchanged values) moves to rows 18-33
#  - Owl block (unchanged values) moves to rows 34-49
$data = @(
    @('Tokay', 0, 'tokay_GG1rearSOAEwf.mat', 1342.7734375, 10.69166059430503, 0.1760964723817636, 0.007962371235322438, 0.0001311438456137643, 0.7695268930878842, 0.0172105369211576, 0.00001986000088344137, 0.024, 32.22656249999999),
    @('Tokay', 0, 'tokay_GG1rearSOAEwf.mat', 1779.1748046875, 17.32916630410013, 0.4232590709113763, 0.009740002083237621, 0.000237896281914648, 0.9570870177067973, 0.03301417938806337, 0.00007663130739317239, 0.03, 53.37524414062499),
    @('Tokay', 0, 'tokay_GG1rearSOAEwf.mat', 3649.90234375, 33.75892247630429, 1.028840650113268, 0.009249267321935947, 0.0002818816924992607, 0.7199080342510847, 0.02708342108361399, 0.00007995199215624267, 0.025, 91.24755859375),
    @('Tokay', 0, 'tokay_GG1rearSOAEwf.mat', 4211.425781249999, 46.12704065170114, 1.075126454008817, 0.01095283237735466, 0.0002552879974272531, 1.12866937685471, 0.03671680144623923, 0.0001037269999680694, 0.032, 134.765625),
    @('Tokay', 1, 'tokay_GG2rearSOAEwf.mat', 1364.1357421875, 11.17892505535753, 0.5165510225718823, 0.008194877320222721, 0.0003786654118039249, 0.6745621062333408, 0.04402790446511812, 0.00007726461382615047, 0.022, 30.01098632812499),
    @('Tokay', 1, 'tokay_GG2rearSOAEwf.mat', 1776.123046875, 15.89647682661615, 0.7054151176192633, 0.008950098842861824, 0.0003971656799681791, 0.7173081035430412, 0.04246533504862162, 0.0001014996746024675, 0.023, 40.85083007812499),
    @('Tokay', 1, 'tokay_GG2rearSOAEwf.mat', 3607.177734375, 37.73114265379662, 1.449612615928588, 0.01046001761827079, 0.0004018689187711335, 0.9157383972568074, 0.03712832976532102, 0.0003092315633094102, 0.025, 90.179443359375),
    @('Tokay', 1, 'tokay_GG2rearSOAEwf.mat', 4394.531249999999, 40.65544364344527, 1.63784606723683, 0.00925137206464177, 0.0003727009717445588, 0.6583620360973251, 0.0297556066450585, 0.0001351758422769943, 0.023, 101.07421875),
    @('Tokay', 2, 'tokay_GG3rearSOAEwf.mat', 1257.32421875, 14.32494353004988, 0.5390601189704403, 0.01139319780564744, 0.0004287359703500823, 0.8770891908533835, 0.04174508411430886, 0.0001964491379633126, 0.03, 37.71972656249999),
    @('Tokay', 2, 'tokay_GG3rearSOAEwf.mat', 1837.158203125, 13.69192523662817, 0.6034139080106731, 0.007452774188601861, 0.0003284496169052116, 0.706938380258768, 0.04688690934595412, 0.00005914898884343351, 0.021, 38.580322265625),
    @('Tokay', 2, 'tokay_GG3rearSOAEwf.mat', 2578.7353515625, 18.99626220687842, 0.6072719896218819, 0.007366503195207008, 0.0002354921722595246, 0.8303156667294969, 0.04097745306654904, 0.00005044628888431499, 0.023, 59.31091308593749),
    @('Tokay', 2, 'tokay_GG3rearSOAEwf.mat', 3567.5048828125, 34.97447838578351, 1.249344345561586, 0.009803624531611243, 0.0003502011592417628, 0.8817297569450167, 0.03434207674933716, 0.0002175897622740973, 0.024, 85.62011718749999),
    @('Tokay', 3, 'tokay_GG4rearSOAEwf.mat', 1251.220703125, 13.65795908537932, 0.5145027674971088, 0.0109157073977978, 0.0004112006508620796, 1.002510682548067, 0.04889588957124038, 0.0002279447290236955, 0.029, 36.28540039062499),
    @('Tokay', 3, 'tokay_GG4rearSOAEwf.mat', 2590.9423828125, 24.41431871821503, 0.8301093520243699, 0.009422949302219909, 0.0003203889663973446, 0.8328716851602599, 0.03451331485304102, 0.0001365347929539875, 0.025, 64.77355957031249),
    @('Tokay', 3, 'tokay_GG4rearSOAEwf.mat', 3216.552734375, 29.62648422299942, 1.131505733697491, 0.009210632210808778, 0.0003517759002068254, 0.6740010612060136, 0.03159363569920032, 0.0001020011721056376, 0.024, 77.19726562499999),
    @('Tokay', 3, 'tokay_GG4rearSOAEwf.mat', 3582.763671875, 34.52374034953, 1.415655296318656, 0.009636064086655871, 0.0003951294101343249, 0.7348881400323389, 0.02998586396090274, 0.0002236340154518327, 0.022, 78.82080078124999),
    @('Anole', 0, 'AC6rearSOAEwfB1.mat', 1235.467529296875, 18.30668826242178, 0.3690403861639351, 0.01481761991174338, 0.0002987050468043964, 1.244989836203278, 0.02952978345395024, 0.0002007714665117125, 0.03900226757369615, 48.18603515625),
    @('Anole', 0, 'AC6rearSOAEwfB1.mat', 2153.3203125, 25.29690995494933, 0.5953484429533717, 0.01174786203803543, 0.0002764792769089581, 1.167710898172541, 0.03010231031943876, 0.0002225317243756103, 0.03, 64.599609375),
    @('Anole', 0, 'AC6rearSOAEwfB1.mat', 3703.7109375, 41.5281063900438, 1.08915715665557, 0.01121256682576724, 0.0002940718579379065, 1.144163426792343, 0.02870813548916876, 0.0003427227295027247, 0.02700680272108844, 100.025390625),
    @('Anole', 0, 'AC6rearSOAEwfB1.mat', 4500.439453125, 41.23098553500826, 0.8730723387695426, 0.009161546547721784, 0.0001939971302498696, 0.8167438705775404, 0.02165438031769892, 0.00005204528376960065, 0.02600907029478458, 117.05224609375),
    @('Anole', 1, 'ACsb4rearSOAEwf1.mat', 966.302490234375, 11.16111458848785, 0.2845045168068702, 0.01155033201433719, 0.000294425937718389, 1.03959466091058, 0.03580911673450345, 0.000121007859575912, 0.03299319727891156, 31.88140869140625),
    @('Anole', 1, 'ACsb4rearSOAEwf1.mat', 3022.723388671875, 32.85410002839859, 0.991437602134656, 0.0108690395394843, 0.0003279948161483192, 1.062929278007326, 0.05328801022286783, 0.00008314328915240177, 0.03399092970521542, 102.7451782226562),
    @('Anole', 1, 'ACsb4rearSOAEwf1.mat', 3151.922607421875, 30.40989913739679, 0.8302988362606313, 0.009648047533207253, 0.0002634261495842301, 1.094492102218201, 0.03294355000433496, 0.0001923642166464746, 0.02399092970521542, 75.6175537109375),
    @('Anole', 1, 'ACsb4rearSOAEwf1.mat', 3954.034423828125, 31.86512210603134, 0.5955162417458386, 0.008058888388528725, 0.0001506097767275596, 0.8485684437948837, 0.01568742471093194, 0.00006257251609513855, 0.02, 79.0806884765625),
    @('Anole', 2, 'ACsb24rearSOAEwfA1.mat', 1811.480712890625, 23.97310956875118, 0.2823202957666863, 0.01323398554461929, 0.0001558505667533058, 0.8390822319957771, 0.01603815453141878, 0.0000120314767802092, 0.04299319727891156, 77.88134765625),
    @('Anole', 2, 'ACsb24rearSOAEwfA1.mat', 2177.545166015625, 28.52641489875884, 0.8644506715203875, 0.0131002632431984, 0.0003969840373516205, 1.259398907356392, 0.03694001935170459, 0.0006131002388727199, 0.03099773242630385, 67.49896240234375),
    @('Anole', 2, 'ACsb24rearSOAEwfA1.mat', 3111.5478515625, 29.68215229911629, 0.4267803169559746, 0.009539352667904833, 0.0001371601329356584, 1.056193478412823, 0.01554751289709087, 0.00006473141480524591, 0.02399092970521542, 74.64892578125),
    @('Anole', 2, 'ACsb24rearSOAEwfA1.mat', 3477.6123046875, 31.06055555179859, 0.9146717045149659, 0.008931575124096447, 0.0002630171578591648, 1.116934099595815, 0.03814322740928681, 0.0002135071431182884, 0.02399092970521542, 83.43115234375),
    @('Anole', 3, 'ACsb30learSOAEwfA2.mat', 1798.0224609375, 19.97564891497187, 0.6034139717733831, 0.01110978830851559, 0.0003355986840446694, 0.974712081647007, 0.03018685437935022, 0.0002432732121475753, 0.02600907029478458, 46.764892578125),
    @('Anole', 3, 'ACsb30learSOAEwfA2.mat', 2139.862060546875, 16.10257133355928, 0.4154622733024121, 0.007525051091117537, 0.0001941537639095458, 0.7013784437903158, 0.02065205830061842, 0.0000543952414301807, 0.02, 42.7972412109375),
    @('Anole', 3, 'ACsb30learSOAEwfA2.mat', 2417.10205078125, 20.27253682937761, 0.5366492113889888, 0.008387124913830249, 0.0002220217434408838, 0.5896820563668638, 0.01678981563699513, 0.00004925455955194371, 0.02099773242630385, 50.753662109375),
    @('Anole', 3, 'ACsb30learSOAEwfA2.mat', 2783.16650390625, 25.43828047953511, 1.144779192487545, 0.009140049811548031, 0.0004113225676152744, 0.5115609906061731, 0.02370136063912024, 0.000125193771606749, 0.02199546485260771, 61.217041015625),
    @('Owl', 0, 'Owl2R1.mat', 4350.5859375, 38.18817150603999, 2.062741876737991, 0.008777707659300775, 0.0004741296704331543, 0.4666897477044611, 0.02632494419016262, 0.0001339445586304853, 0.021, 91.3623046875),
    @('Owl', 0, 'Owl2R1.mat', 7453.125, 69.19801128502411, 1.29562100210892, 0.009284429187089188, 0.0001738359415827481, 0.5678163655099446, 0.01095659896956142, 0.00002903193464227278, 0.023, 171.421875),
    @('Owl', 0, 'Owl2R1.mat', 8452.1484375, 69.4790088667789, 1.153285808860182, 0.008220277883256105, 0.000136448835155728, 0.6154222216427551, 0.008877207013722566, 0.00003359224672051262, 0.019, 160.5908203125),
    @('Owl', 0, 'Owl2R1.mat', 9026.3671875, 81.38665497811712, 2.179639449356527, 0.00901654600211965, 0.0002414747155838022, 0.721227100608454, 0.01785733306857111, 0.0001004126253638873, 0.02, 180.52734375),
    @('Owl', 1, 'Owl7L1.mat', 6837.890625, 61.14423462696202, 2.43283311556865, 0.008941973187376337, 0.000355787076598443, 0.5542966708632717, 0.02302424109636626, 0.0001126053621516931, 0.022, 150.43359375),
    @('Owl', 1, 'Owl7L1.mat', 7901.3671875, 72.61693277193591, 2.621698804344461, 0.009190426295694274, 0.0003318031857185426, 0.7176820945998618, 0.02887872450328945, 0.0001181910140158244, 0.022, 173.830078125),
    @('Owl', 1, 'Owl7L1.mat', 8835.9375, 75.33983965195233, 2.557140366460105, 0.008526524735145799, 0.0002894022695905335, 0.7128489007584123, 0.02314669123347024, 0.0001519252691792823, 0.02, 176.71875),
    @('Owl', 1, 'Owl7L1.mat', 9257.8125, 89.1534077998639, 3.259690495762425, 0.009630072741251121, 0.0003521015894156881, 0.3687694104084971, 0.01107179803851206, 0.00007302487491588827, 0.021, 194.4140625),
    @('Owl', 2, 'TAG6rearSOAEwf1.mat', 5625.54931640625, 51.83701769448351, 2.526739557456056, 0.009214569952005749, 0.0004491542808250064, 0.7835877415703257, 0.04687451605713499, 0.0002243659335072323, 0.02399092970521542, 134.962158203125),
    @('Owl', 2, 'TAG6rearSOAEwf1.mat', 8096.484375, 78.28887462547326, 2.624205836221664, 0.009669489990891665, 0.0003241167048163005, 0.606666242649432, 0.02019008849351071, 0.0001021461108885125, 0.02199546485260771, 178.0859375),
    @('Owl', 2, 'TAG6rearSOAEwf1.mat', 8489.46533203125, 101.9184343584824, 3.481792421635467, 0.01200528306228405, 0.0004101309429344697, 0.941342888221613, 0.03199093798015184, 0.0003485549518335093, 0.02800453514739229, 237.7435302734375),
    @('Owl', 2, 'TAG6rearSOAEwf1.mat', 9864.898681640625, 92.57374287449112, 3.26948592510969, 0.009384155464949513, 0.0003314262042239185, 0.6031090155337008, 0.02156217105194045, 0.0001095498729547878, 0.02199546485260771, 216.9830322265625),
    @('Owl', 3, 'owl_TAG4learSOAEwf1.mat', 4944.561767578125, 45.28841458447616, 1.798293050157823, 0.009159237302168175, 0.0003636910882475713, 0.5391445149782679, 0.02201095677423928, 0.0001084582128997591, 0.02199546485260771, 108.7579345703125),
    @('Owl', 3, 'owl_TAG4learSOAEwf1.mat', 5768.206787109375, 68.8618131122638, 2.105416667748536, 0.01193816651409832, 0.0003650036736639992, 0.9243709592246925, 0.03284198698865044, 0.0002025271991541933, 0.03099773242630385, 178.8013305664062),
    @('Owl', 3, 'owl_TAG4learSOAEwf1.mat', 7184.014892578125, 215.212239384061, 3.426336935219513, 0.02995709816893599, 0.0004769390078463361, 1.272598354498129, 0.01906358170377451, 0.0004631125362880487, 0.07099773242630386, 510.0487670898438),
    @('Owl', 3, 'owl_TAG4learSOAEwf1.mat', 9633.416748046875, 163.3124452650469, 3.053847747422519, 0.01695270219656568, 0.0003170056717458705, 1.103992840979661, 0.01949410344315334, 0.0002697444831906281, 0.04099773242630386, 394.9482421875)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
